{"js": "// Replace the two-digit multiplication equations in the table cells with\n// their new values, per the commit's regenerated answer set.\nconst replacements = [\n  [\"14\u00d782=1148\", \"15\u00d715=225\"],\n  [\"55\u00d719=1045\", \"90\u00d743=3870\"],\n  [\"93\u00d796=8928\", \"24\u00d799=2376\"],\n  [\"50\u00d741=2050\", \"23\u00d743=989\"],\n  [\"13\u00d773=949\", \"50\u00d794=4700\"],\n  [\"64\u00d717=1088\", \"24\u00d723=552\"],\n  [\"63\u00d732=2016\", \"77\u00d759=4543\"],\n  [\"25\u00d782=2050\", \"97\u00d779=7663\"],\n  [\"48\u00d747=2256\", \"86\u00d733=2838\"],\n  [\"56\u00d732=1792\", \"33\u00d739=1287\"],\n  [\"88\u00d726=2288\", \"48\u00d741=1968\"],\n  [\"30\u00d712=360\", \"13\u00d753=689\"],\n  [\"94\u00d781=7614\", \"96\u00d753=5088\"],\n  [\"65\u00d746=2990\", \"14\u00d777=1078\"],\n  [\"35\u00d754=1890\", \"13\u00d726=338\"],\n  [\"73\u00d740=2920\", \"31\u00d750=1550\"],\n  [\"99\u00d787=8613\", \"36\u00d716=576\"],\n  [\"30\u00d792=2760\", \"85\u00d754=4590\"],\n  [\"36\u00d778=2808\", \"44\u00d765=2860\"],\n  [\"50\u00d768=3400\", \"24\u00d756=1344\"],\n  [\"64\u00d763=4032\", \"88\u00d711=968\"],\n  [\"95\u00d722=2090\", \"43\u00d763=2709\"],\n  [\"21\u00d765=1365\", \"63\u00d793=5859\"],\n  [\"64\u00d739=2496\", \"91\u00d733=3003\"],\n  [\"25\u00d758=1450\", \"92\u00d783=7636\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication equations in the table cells with\n# their new values, per the commit's regenerated answer set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"14\u00d782=1148\", \"15\u00d715=225\"),\n  @(\"55\u00d719=1045\", \"90\u00d743=3870\"),\n  @(\"93\u00d796=8928\", \"24\u00d799=2376\"),\n  @(\"50\u00d741=2050\", \"23\u00d743=989\"),\n  @(\"13\u00d773=949\",  \"50\u00d794=4700\"),\n  @(\"64\u00d717=1088\", \"24\u00d723=552\"),\n  @(\"63\u00d732=2016\", \"77\u00d759=4543\"),\n  @(\"25\u00d782=2050\", \"97\u00d779=7663\"),\n  @(\"48\u00d747=2256\", \"86\u00d733=2838\"),\n  @(\"56\u00d732=1792\", \"33\u00d739=1287\"),\n  @(\"88\u00d726=2288\", \"48\u00d741=1968\"),\n  @(\"30\u00d712=360\",  \"13\u00d753=689\"),\n  @(\"94\u00d781=7614\", \"96\u00d753=5088\"),\n  @(\"65\u00d746=2990\", \"14\u00d777=1078\"),\n  @(\"35\u00d754=1890\", \"13\u00d726=338\"),\n  @(\"73\u00d740=2920\", \"31\u00d750=1550\"),\n  @(\"99\u00d787=8613\", \"36\u00d716=576\"),\n  @(\"30\u00d792=2760\", \"85\u00d754=4590\"),\n  @(\"36\u00d778=2808\", \"44\u00d765=2860\"),\n  @(\"50\u00d768=3400\", \"24\u00d756=1344\"),\n  @(\"64\u00d763=4032\", \"88\u00d711=968\"),\n  @(\"95\u00d722=2090\", \"43\u00d763=2709\"),\n  @(\"21\u00d765=1365\", \"63\u00d793=5859\"),\n  @(\"64\u00d739=2496\", \"91\u00d733=3003\"),\n  @(\"25\u00d758=1450\", \"92\u00d783=7636\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
